# CaUWMET_Tests.xlsx edit script
# - Rename the two existing sheets
# - Add a new "TestInputData" sheet with the raw test inputs
# - Rework the "test_modelLogic.py" sheet so its values/formulas pull from
#   the new TestInputData sheet
# - Restore view/selection state to match the target workbook

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "test_storageOperations.py"
$ws2.Name = "test_modelLogic.py"

# --- New sheet: TestInputData -------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "TestInputData"

$ws3.Range("A1").Value = "i"
$ws3.Range("B1").Value = 0

$ws3.Range("A2").Value = "Hydrologic Year Type at i"
$ws3.Range("B2").Value = "W"

$ws3.Range("A3").Value = "Normal or Better Demands (acre-feet/year)"
$ws3.Range("B3").Value = 1000000

$ws3.Range("A4").Value = "Base Long-term Conservation (acre-feet/year)"
$ws3.Range("B4").Value = 5

$ws3.Range("A5").Value = "Base Local Supply (Total, acre-feet/year)"
$ws3.Range("B5").Formula = "=7*5000"

$ws3.Range("A6").Value = "SWP/CVP Supply"
$ws3.Range("B6").Value = 5000

$ws3.Range("A8").Value = "Surface initial storage (acre-feet)"
$ws3.Range("B8").Value = 516500

$ws3.Range("A9").Value = "Groundwater initial storage (acre-feet)"
$ws3.Range("B9").Value = 95000

$ws3.Range("A10").Value = "Surface max take capacity (acre-feet)"
$ws3.Range("B10").Value = 815000

$ws3.Range("A11").Value = "Groundwater max take capacity (acre-feet)"
$ws3.Range("B11").Value = 194300

$ws3.Columns.Item(1).ColumnWidth = 42.7109375

# --- Rework test_modelLogic.py -------------------------------------------
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "Variable"
$ws2.Range("B1").Value = "Value Used In Test"

$ws2.Range("A2").Value = "appliedDemand"
$ws2.Range("B2").Formula = "=TestInputData!B3-TestInputData!B4"

$ws2.Range("A3").Value = "demandsToBeMetBySWPCVP"
$ws2.Range("B3").Formula = "=B2-TestInputData!B5"

$ws2.Range("A4").Value = "demandsToBeMetByStorage"
$ws2.Range("B4").Formula = "=B3-TestInputData!B6"

$ws2.Range("A5").Value = "demandsToBeMetByContingentOptions"
$ws2.Range("B5").Formula = "=MAX(0, B4-MIN(TestInputData!B8,TestInputData!B10) - MIN(TestInputData!B9,TestInputData!B11))"

$ws2.Range("A12").Value = "pctCapacitySurfaceCarryover_Contractor"
$ws2.Range("B12").Value = 0.1

$ws2.Range("A13").Value = "pctStorageCalledSurfaceCarryover_Contractor"
$ws2.Range("B13").Value = 50

$ws2.Range("A15").Value = "pctCapacityGroundwaterBank_Contractor"
$ws2.Range("B15").Formula = "=B2*B12"

# --- View / selection state ----------------------------------------------
$ws3.Range("A13").Select()
$ws2.Range("A10").Select()
$ws1.Activate()
$ws1.Range("C15").Select()

$wb.Application.Calculate()
